$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes a value as plain text, preventing Excel's automatic number/date
# coercion for values that look numeric (e.g. "242.34", "5.239"), and then
# resets the cell's style back to Normal/General so no stray formatting
# is left behind on the cell.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row -> (new Price [col D], new Volume(1h) [col E]); $null means "no change"
$changes = @{
    2  = @("30.008.34", "  -1.02%  ")
    3  = @("1.900.42", "  -1.44%  ")
    4  = @($null, "  -0.48%  ")
    5  = @("0.7425", "  +0.99%  ")
    6  = @("242.34", "  -0.50%  ")
    7  = @($null, "  -0.25%  ")
    8  = @("0.3069", "  -2.37%  ")
    9  = @($null, "  -7.07%  ")
    10 = @("0.06892", "  -1.61%  ")
    11 = @("0.08017", "  -0.21%  ")
    12 = @("0.7534", "  -2.59%  ")
    13 = @("1.901.65", "  -1.13%  ")
    14 = @("5.239", "  -2.16%  ")
    15 = @("91.17", "  -1.86%  ")
    16 = @("6.173", "  +3.91%  ")
    17 = @("30.024.73", "  -0.97%  ")
    18 = @("14.01", "  -2.75%  ")
    19 = @("0.000007768", "  -1.82%  ")
    20 = @("237.10", "  -5.09%  ")
    22 = @("2.157.75", "  +0.08%  ")
    23 = @($null, "  -0.37%  ")
    24 = @("7.096", "  +7.00%  ")
    25 = @("9.345", "  -2.00%  ")
    26 = @("167.64", "  +1.23%  ")
    27 = @("18.79", "  -1.01%  ")
    28 = @("0.1269", "  -0.56%  ")
    29 = @("2.042", "  -5.50%  ")
    30 = @($null, "  -0.57%  ")
    31 = @($null, "  -1.45%  ")
    32 = @("4.299", "  -2.35%  ")
    33 = @("4.059", "  -0.99%  ")
    34 = @("0.05344", "  +2.72%  ")
    35 = @("1.283", "  -1.58%  ")
    36 = @("0.7379", "  -1.87%  ")
    37 = @("2.726", "  -1.62%  ")
    38 = @("0.01946", "  +0.12%  ")
    39 = @("2.764", "  -0.95%  ")
    40 = @("6.256", "  -3.88%  ")
    41 = @("0.4458", "  -0.57%  ")
    42 = @("72.74", "  -4.80%  ")
    43 = @("1.947", "  +0.06%  ")
    44 = @($null, "  -0.23%  ")
    45 = @("7.764", "  +1.36%  ")
    46 = @("0.8319", "  -0.77%  ")
    47 = @("101.61", "  +0.37%  ")
    48 = @("9.845", "  -0.90%  ")
    49 = @("2.052.92", "  -0.75%  ")
    50 = @("36.59", "  -2.21%  ")
    51 = @("0.05987", "  -0.52%  ")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $priceVal = $vals[0]
    $volVal = $vals[1]

    if ($null -ne $priceVal) {
        Set-TextValue $ws.Cells.Item($row, 4) $priceVal
    }
    if ($null -ne $volVal) {
        Set-TextValue $ws.Cells.Item($row, 5) $volVal
    }
}
